$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lecture Attendance Tracker")
$full = $ws.Range("C2:AG11")
$fcs = $full.FormatConditions
for ($i = 1; $i -le $fcs.Count(); $i++) {
  $fc = $fcs.Item($i)
  $t = $fc.Type()
  $op = $fc.Operator()
  Write-Output "idx=$i type=$t op=$op"
}
